# Update the "Förändrad" (Changed) date column (C) for all existing data
# rows (2-233) from 45205 to 45206.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C233").Value = 45206

# Row 233 picks up an explicit row height (15, custom) in the edited file.
$ws.Rows.Item(233).RowHeight = 15

# Copy the date-format / wrap-text formatting from row 233 down to the new
# row 234 before filling in its values, so B234/C234 keep the date style
# and R234 keeps the wrap-text style (matching an empty placeholder cell).
$ws.Range("B233:C233").Copy()
$ws.Range("B234").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("R233").Copy()
$ws.Range("R234").PasteSpecial(-4122)  # xlPasteFormats

# Append the new record as the last row of the table.
$ws.Range("A234").Value = "A 47978-2023"
$ws.Range("B234").Value = 45204
$ws.Range("C234").Value = 45206
$ws.Range("D234").Value = "ÖSTERGÖTLANDS LÄN"
$ws.Range("E234").Value = "BOXHOLM"
$ws.Range("G234").Value = 1
$ws.Range("H234").Value = 0
$ws.Range("I234").Value = 0
$ws.Range("J234").Value = 0
$ws.Range("K234").Value = 0
$ws.Range("L234").Value = 0
$ws.Range("M234").Value = 0
$ws.Range("N234").Value = 0
$ws.Range("O234").Value = 0
$ws.Range("P234").Value = 0
$ws.Range("Q234").Value = 0
